$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 49

# Text columns that look numeric must be forced to text so Excel
# doesn't silently coerce them to Number/Date types.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "6376"

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "7/8/2025"

$ws.Cells.Item($row, 3).Value = "BOYACA 712"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "7"

$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Cells.Item($row, 5).Value = "808099366"

$ws.Cells.Item($row, 6).Value = "PEBCOM"
$ws.Cells.Item($row, 7).Value = "Pendiente"
$ws.Cells.Item($row, 8).Value = "Picada"

$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 9).Value = "1"

$ws.Cells.Item($row, 10).Value = "Cambio"
$ws.Cells.Item($row, 11).Value = "Sin equipos"
$ws.Cells.Item($row, 12).Value = "Pasante"

$ws.Cells.Item($row, 13).Value = -58.461858
$ws.Cells.Item($row, 14).Value = -34.619348

$ws.Cells.Item($row, 15).Value = "Boedo"
$ws.Cells.Item($row, 16).Value = "Capital Sur"
